$d = $word.ActiveDocument

# Locate the title paragraph ("DSE 511 Final Project Report") so the new
# team-name paragraph can be inserted directly after it, regardless of
# exact paragraph indices.
$titleRange = $d.Content
$found = $titleRange.Find.Execute("DSE 511 Final Project Report", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find the title paragraph to anchor the new heading on."
}

# Insert a brand-new (empty) paragraph right after the title.
$titleRange.InsertParagraphAfter()

# That new paragraph is now the second paragraph in the document; grab its
# Range so we can populate it with the formatted "Lancaster Barnstormers:"
# heading (three runs, matching the target markup exactly, including the
# middle run's eastAsia font hint).
$newPara = $d.Paragraphs(2)
$newRange = $newPara.Range

$xml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:rFonts w:cs="Times New Roman"/><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cs="Times New Roman"/><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>L</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="Times New Roman" w:hint="eastAsia"/><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>an</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="Times New Roman"/><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>caster Barnstormers:</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$newRange.InsertXML($xml)
